# Refresh the cryptocurrency table (columns D "Price" and E "Volume(1h)")
# with the latest scrape, per the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.063.97"
$ws.Range("E2").Value = "  -0.45%  "

$ws.Range("D3").Value = "3.664.26"
$ws.Range("E3").Value = "  -1.28%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "2.42"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +28.19%  "

$ws.Range("E5").Value = "  +0.16%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "227.83"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.63%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "646.41"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.68%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.430"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.59%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "1.13"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +7.20%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("D11").Value = "3.661.30"
$ws.Range("E11").Value = "  -1.36%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "47.52"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +7.21%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.210"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.84%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000299"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.16%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.62"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.32%  "

$ws.Range("D16").Value = "4.363.46"
$ws.Range("E16").Value = "  -0.78%  "

$ws.Range("D17").Value = "95.937.48"
$ws.Range("E17").Value = "  -0.30%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "8.87"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").Value = "3.665.20"
$ws.Range("E19").Value = "  -0.64%  "

$ws.Range("E20").Value = "  +5.58%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.92"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.23%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.535"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +6.07%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "523.89"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.59%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.30"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.86%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.248"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +39.64%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "121.79"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +20.23%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0000207"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.37%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "6.82"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.32%  "

$ws.Range("D29").Value = "3.862.11"
$ws.Range("E29").Value = "  -1.07%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "12.98"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.68%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "13.19"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +8.89%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.99"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.39%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.184"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.87%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.82"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.86%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "32.79"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.48%  "

# Row 37: coin identity swapped with its neighbour, then price refreshed
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.20%  "

# Row 38: coin identity swapped with its neighbour, then price refreshed
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.611"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.75%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "613.27"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.90%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "8.47"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.69%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "7.06"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.16%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.498"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +15.45%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.163"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.09%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0502"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +11.09%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "39.97"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.12%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.99"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.01%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.954"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.67%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.99"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +5.93%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.27"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("E51").Value = "  -0.26%  "
